# Update cryptos list values (price and volume-1h columns) per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.718.50"
$ws.Range("E2").Value = "  -0.55%  "

$ws.Range("D3").Value = "2.581.21"
$ws.Range("E3").Value = "  +1.17%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'581.53"
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("D6").Value = "'144.85"
$ws.Range("E6").Value = "  -1.50%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "'0.593"
$ws.Range("E8").Value = "  +1.46%  "

$ws.Range("E9").Value = "  +0.37%  "

$ws.Range("D10").Value = "'5.55"
$ws.Range("E10").Value = "  -0.46%  "

$ws.Range("D12").Value = "'0.352"
$ws.Range("E12").Value = "  -0.40%  "

$ws.Range("D13").Value = "'26.93"
$ws.Range("E13").Value = "  -2.22%  "

$ws.Range("D14").Value = "3.043.60"
$ws.Range("E14").Value = "  +1.22%  "

$ws.Range("D15").Value = "62.595.20"
$ws.Range("E15").Value = "  -0.60%  "

$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("D17").Value = "2.579.07"
$ws.Range("E17").Value = "  +1.23%  "

$ws.Range("D18").Value = "'11.22"
$ws.Range("E18").Value = "  -1.09%  "

$ws.Range("D19").Value = "'338.37"
$ws.Range("E19").Value = "  -0.18%  "

$ws.Range("D20").Value = "'4.35"
$ws.Range("E20").Value = "  +0.70%  "

$ws.Range("D21").Value = "'6.65"
$ws.Range("E21").Value = "  -1.57%  "

$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").Value = "'67.12"
$ws.Range("E23").Value = "  +2.45%  "

$ws.Range("D24").Value = "2.701.33"
$ws.Range("E24").Value = "  +0.93%  "

$ws.Range("E25").Value = "  -1.88%  "

$ws.Range("E26").Value = "  -2.45%  "

$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("E28").Value = "  -0.69%  "

$ws.Range("D29").Value = "'7.92"
$ws.Range("E29").Value = "  +2.89%  "

$ws.Range("E30").Value = "  -1.81%  "

$ws.Range("E31").Value = "  -2.67%  "

$ws.Range("D32").Value = "0.0₃0809"
$ws.Range("E32").Value = "  -0.77%  "

$ws.Range("D33").Value = "'461.03"
$ws.Range("E33").Value = "  +9.30%  "

$ws.Range("D34").Value = "'176.59"
$ws.Range("E34").Value = "  -0.79%  "

$ws.Range("D35").Value = "'1.62"
$ws.Range("E35").Value = "  +4.07%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").Value = "'0.401"
$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("D38").Value = "'18.89"
$ws.Range("E38").Value = "  -1.05%  "

$ws.Range("D39").Value = "'4.45"
$ws.Range("E39").Value = "  +1.85%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("E41").Value = "  -3.16%  "

$ws.Range("D42").Value = "'157.26"
$ws.Range("E42").Value = "  +4.27%  "

$ws.Range("D43").Value = "'3.72"
$ws.Range("E43").Value = "  -1.55%  "

$ws.Range("D44").Value = "'21.12"
$ws.Range("E44").Value = "  +1.70%  "

$ws.Range("E45").Value = "  +3.80%  "

$ws.Range("D46").Value = "'0.0536"
$ws.Range("E46").Value = "  -0.54%  "

$ws.Range("E47").Value = "  -0.37%  "

$ws.Range("E48").Value = "  -1.94%  "

$ws.Range("D49").Value = "'18.10"
$ws.Range("E49").Value = "  -1.11%  "

$ws.Range("D50").Value = "'11.41"
$ws.Range("E50").Value = "  +0.86%  "

$ws.Range("E51").Value = "  -1.00%  "
